$wb = $excel.ActiveWorkbook

$wsNov = $wb.Worksheets.Item("NOV-2020")
$wsDec = $wb.Worksheets.Item("DEC-2020")

# --- Populate the two new daily rows (Dec 1 & Dec 2) on the DEC-2020 sheet ---
# Row 10 of NOV-2020 already carries the exact A/B/C/D/E/F/G "completed row"
# style pattern these two new rows need, so copy its formatting over first.
$wsNov.Range("A10:G10").Copy() | Out-Null
$wsDec.Range("A2:G2").PasteSpecial(-4122) | Out-Null
$wsNov.Range("A10:G10").Copy() | Out-Null
$wsDec.Range("A3:G3").PasteSpecial(-4122) | Out-Null

$wsDec.Range("A2").Value = 1
$wsDec.Range("B2").Value = 44166
$wsDec.Range("C2").Value = "SONIYA, nQGCare , nMVAR"
$wsDec.Range("D2").Value = "SONIYA setup tested. nQGCare Setup files created"
$wsDec.Range("E2").Value = 1
$wsDec.Range("F2").Value = "COMPLETED"

$wsDec.Range("A3").Value = 2
$wsDec.Range("B3").Value = 44167
$wsDec.Range("C3").Value = " nQGCare , nMVAR"
$wsDec.Range("D3").Value = "nQGCare,nMVAR Setup files created"
$wsDec.Range("E3").Value = 1
$wsDec.Range("F3").Value = "COMPLETED"

# --- View / selection state ---
# NOV-2020 is scrolled further down and the selection moves to A31:G31; it is
# no longer the active tab.
$wsNov.Activate() | Out-Null
$wsNov.Range("A22").Select() | Out-Null
$wsNov.Range("A31:G31").Select() | Out-Null

# DEC-2020 becomes the active tab with D3 selected.
$wsDec.Activate() | Out-Null
$wsDec.Range("D3").Select() | Out-Null
